# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right before "2022-Q3" (duplicating the
#    "2022-Q3" sheet so formatting/styles carry over, then overwriting the
#    data rows with the 2022-Q4 figures, and adding two extra data rows).
# 2. Update the "总计" (summary) sheet: insert the 2022-Q4 summary row at the
#    top of the data (row 2), pushing the older quarters down by one row and
#    re-numbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: summary sheet ("总计") — shift existing rows down and add 2022-Q4
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Final target state (already includes the new 2022-Q4 row at the top and the
# 2020-Q4 row moved to row 7):
#   row2: 0  2022-Q4  7  1.91
#   row3: 1  2022-Q3  5  0.13
#   row4: 2  2021-Q3  1  0
#   row5: 3  2021-Q2  1  0
#   row6: 4  2021-Q1  2  0.01
#   row7: 5  2020-Q4  2  0.05

# Row 7 is brand new — give column A the same style as the existing index
# cells (bold/centered/bordered, cellXfs index "2") before writing the value.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 0.05

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.01

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 0.13

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 1.91

# ---------------------------------------------------------------------------
# Part 2: new "2022-Q4" worksheet, inserted before "2022-Q3"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

# Duplicate the 2022-Q3 sheet so the header row / column-A styling (and the
# general per-cell formatting) match exactly; the copy is placed immediately
# before the source sheet.
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Two more data rows are needed (8 total instead of 6) — extend column A's
# styling down first.
$q4.Range("A6").Copy()
$q4.Range("A7").PasteSpecial(-4122)
$q4.Range("A6").Copy()
$q4.Range("A8").PasteSpecial(-4122)

function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2: 010591 富国中国中小盘混合（QDII）美元
$q4.Range("A2").Value = 0
Set-TextCell $q4.Range("B2") "010591"
$q4.Range("C2").Value = "富国中国中小盘混合（QDII）美元"
Set-TextCell $q4.Range("D2") "32.21"
Set-TextCell $q4.Range("E2") "87.21"
Set-TextCell $q4.Range("F2") "2.78"
Set-TextCell $q4.Range("G2") "0.8954"
$q4.Range("H2").Value = 5

# Row 3: 100061 富国中国中小盘混合（QDII）人民币
$q4.Range("A3").Value = 1
Set-TextCell $q4.Range("B3") "100061"
$q4.Range("C3").Value = "富国中国中小盘混合（QDII）人民币"
Set-TextCell $q4.Range("D3") "32.21"
Set-TextCell $q4.Range("E3") "87.21"
Set-TextCell $q4.Range("F3") "2.78"
Set-TextCell $q4.Range("G3") "0.8954"
$q4.Range("H3").Value = 5

# Row 4: 005051 上投摩根标普港股通低波红利指数A
$q4.Range("A4").Value = 2
Set-TextCell $q4.Range("B4") "005051"
$q4.Range("C4").Value = "上投摩根标普港股通低波红利指数A"
Set-TextCell $q4.Range("D4") "1.64"
Set-TextCell $q4.Range("E4") "93.98"
Set-TextCell $q4.Range("F4") "2.42"
Set-TextCell $q4.Range("G4") "0.0397"
$q4.Range("H4").Value = 7

# Row 5: 005052 上投摩根标普港股通低波红利指数C
$q4.Range("A5").Value = 3
Set-TextCell $q4.Range("B5") "005052"
$q4.Range("C5").Value = "上投摩根标普港股通低波红利指数C"
Set-TextCell $q4.Range("D5") "1.52"
Set-TextCell $q4.Range("E5") "93.98"
Set-TextCell $q4.Range("F5") "2.42"
Set-TextCell $q4.Range("G5") "0.0368"
$q4.Range("H5").Value = 7

# Row 6: 006477 中邮沪港深精选混合
$q4.Range("A6").Value = 4
Set-TextCell $q4.Range("B6") "006477"
$q4.Range("C6").Value = "中邮沪港深精选混合"
Set-TextCell $q4.Range("D6") "0.67"
Set-TextCell $q4.Range("E6") "94.24"
Set-TextCell $q4.Range("F6") "5.20"
Set-TextCell $q4.Range("G6") "0.0348"
$q4.Range("H6").Value = 7

# Row 7: 501303 广发恒生中型股指数（LOF）A
$q4.Range("A7").Value = 5
Set-TextCell $q4.Range("B7") "501303"
$q4.Range("C7").Value = "广发恒生中型股指数（LOF）A"
Set-TextCell $q4.Range("D7") "0.24"
Set-TextCell $q4.Range("E7") "90.77"
Set-TextCell $q4.Range("F7") "1.28"
Set-TextCell $q4.Range("G7") "0.0031"
$q4.Range("H7").Value = 9

# Row 8: 004996 广发恒生中型股指数（LOF）C
$q4.Range("A8").Value = 6
Set-TextCell $q4.Range("B8") "004996"
$q4.Range("C8").Value = "广发恒生中型股指数（LOF）C"
Set-TextCell $q4.Range("D8") "0.12"
Set-TextCell $q4.Range("E8") "90.77"
Set-TextCell $q4.Range("F8") "1.28"
Set-TextCell $q4.Range("G8") "0.0015"
$q4.Range("H8").Value = 9
